# Updated cryptos list - apply cell value changes per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.989.49"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "'2.051.57"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'248.64"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'0.671"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'54.77"
$ws.Range("E8").Value = "  +15.27%  "
$ws.Range("D9").Value = "'60.48"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "'0.380"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +5.08%  "
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "'2.351.39"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "'2.053.95"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "'36.957.46"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'0.0₃0927"
$ws.Range("E19").Value = "  +12.10%  "
$ws.Range("D20").Value = "'72.91"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'14.19"
$ws.Range("E21").Value = "  +8.24%  "
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("D23").Value = "'235.82"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "'169.87"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'8.94"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'19.98"
$ws.Range("E28").Value = "  -6.77%  "
$ws.Range("D29").Value = "'1.98"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("E34").Value = "  +5.90%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.0867"
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("E37").Value = "  -5.37%  "
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").Value = "'0.103"
$ws.Range("E40").Value = "  +22.43%  "
$ws.Range("D41").Value = "'17.72"
$ws.Range("E41").Value = "  +11.05%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "'96.07"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "'2.80"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "'4.19"
$ws.Range("E46").Value = "  +46.78%  "
$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").Value = "'13.24"
$ws.Range("E47").Value = "  -52.91%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  +6.19%  "
$ws.Range("D49").Value = "'1.292.43"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").Value = "'4.08"
$ws.Range("E51").Value = "  +6.57%  "
